$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (ALC)
$ws.Cells.Item(2, 8).Value = 276.46155  # H2: 329.875 -> 276.46155
$ws.Cells.Item(2, 9).Value = 185.7  # I2: 180.4 -> 185.7
$ws.Cells.Item(2, 11).Value = 185.7  # K2: 180.4 -> 185.7
$ws.Cells.Item(2, 13).Value = -72.69999999999999  # M2: -67.40000000000001 -> -72.69999999999999

# Row 6 (ALC)
$ws.Cells.Item(6, 8).Value = 1116.6923  # H6: 1038.1428 -> 1116.6923
$ws.Cells.Item(6, 9).Value = 279.77777  # I6: 253.5 -> 279.77777
$ws.Cells.Item(6, 11).Value = 839.33331  # K6: 760.5 -> 839.33331
$ws.Cells.Item(6, 13).Value = -727.33331  # M6: -648.5 -> -727.33331

# Row 9 (ALC)
$ws.Cells.Item(9, 8).Value = 183.5  # H9: 357.5 -> 183.5
$ws.Cells.Item(9, 9).Value = 218  # I9: 630 -> 218
$ws.Cells.Item(9, 10).Value = 80  # J9: 85 -> 80
$ws.Cells.Item(9, 11).Value = 218  # K9: 630 -> 218
$ws.Cells.Item(9, 12).Value = 80  # L9: 85 -> 80
$ws.Cells.Item(9, 13).Value = -49  # M9: -461 -> -49
$ws.Cells.Item(9, 14).Value = -418  # N9: -423 -> -418

# Row 10 (ALC)
$ws.Cells.Item(10, 8).Value = 0  # H10: 8684 -> 0
$ws.Cells.Item(10, 9).Value = 0  # I10: 5000 -> 0
$ws.Cells.Item(10, 10).Value = 0  # J10: 10526 -> 0
$ws.Cells.Item(10, 11).Value = 0  # K10: 5000 -> 0
$ws.Cells.Item(10, 12).Value = 0  # L10: 10526 -> 0
$ws.Cells.Item(10, 13).ClearContents()  # M10: -4707 -> (removed)
$ws.Cells.Item(10, 14).ClearContents()  # N10: -11112 -> (removed)

# Row 12 (ALC)
$ws.Cells.Item(12, 8).Value = 834  # H12: 522.1667 -> 834
$ws.Cells.Item(12, 9).Value = 500  # I12: 367 -> 500
$ws.Cells.Item(12, 10).Value = 1001  # J12: 677.3333 -> 1001
$ws.Cells.Item(12, 11).Value = 500  # K12: 367 -> 500
$ws.Cells.Item(12, 12).Value = 1001  # L12: 677.3333 -> 1001
$ws.Cells.Item(12, 13).Value = -330  # M12: -197 -> -330
$ws.Cells.Item(12, 14).Value = -1341  # N12: -1017.3333 -> -1341

# Row 13 (ALC)
$ws.Cells.Item(13, 8).Value = 1198.6666  # H13: 34918 -> 1198.6666
$ws.Cells.Item(13, 10).Value = 1198.6666  # J13: 34918 -> 1198.6666
$ws.Cells.Item(13, 12).Value = 1198.6666  # L13: 34918 -> 1198.6666
$ws.Cells.Item(13, 14).Value = -1536.6666  # N13: -35256 -> -1536.6666

# Row 28 (ALC)
$ws.Cells.Item(28, 8).Value = 5833.1665  # H28: 6152.8237 -> 5833.1665
$ws.Cells.Item(28, 9).Value = 957.6667  # I28: 1008.4545 -> 957.6667
$ws.Cells.Item(28, 11).Value = 957.6667  # K28: 1008.4545 -> 957.6667
$ws.Cells.Item(28, 13).Value = -472.6667  # M28: -523.4545000000001 -> -472.6667

# Row 38 (ALC)
$ws.Cells.Item(38, 8).Value = 336.2857  # H38: 373.22223 -> 336.2857
$ws.Cells.Item(38, 9).Value = 336.2857  # I38: 294.875 -> 336.2857
$ws.Cells.Item(38, 10).Value = 0  # J38: 1000 -> 0
$ws.Cells.Item(38, 11).Value = 1008.8571  # K38: 884.625 -> 1008.8571
$ws.Cells.Item(38, 12).Value = 0  # L38: 3000 -> 0
$ws.Cells.Item(38, 13).Value = -636.8571000000001  # M38: -512.625 -> -636.8571000000001
$ws.Cells.Item(38, 14).ClearContents()  # N38: -3744 -> (removed)

# Row 52 (ALC)
$ws.Cells.Item(52, 8).Value = 397.4  # H52: 398.125 -> 397.4
$ws.Cells.Item(52, 9).Value = 394.5  # I52: 0 -> 394.5
$ws.Cells.Item(52, 11).Value = 1183.5  # K52: 0 -> 1183.5
$ws.Cells.Item(52, 13).Value = -1023.5  # M52: None -> -1023.5

# Row 58 (ALC)
$ws.Cells.Item(58, 8).Value = 1729.5555  # H58: 1856.6 -> 1729.5555
$ws.Cells.Item(58, 10).Value = 2999.8  # J58: 2999.8333 -> 2999.8
$ws.Cells.Item(58, 12).Value = 8999.400000000001  # L58: 8999.499899999999 -> 8999.400000000001
$ws.Cells.Item(58, 14).Value = -9299.400000000001  # N58: -9299.499899999999 -> -9299.400000000001

# Row 61 (ALC)
$ws.Cells.Item(61, 8).Value = 14172  # H61: 35207.5 -> 14172
$ws.Cells.Item(61, 9).Value = 14172  # I61: 35207.5 -> 14172
$ws.Cells.Item(61, 11).Value = 42516  # K61: 105622.5 -> 42516
$ws.Cells.Item(61, 13).Value = -42344  # M61: -105450.5 -> -42344

# Row 103 (ALC)
$ws.Cells.Item(103, 8).Value = 1499.6666  # H103: 3214.1428 -> 1499.6666
$ws.Cells.Item(103, 9).Value = 1000  # I103: 3750 -> 1000
$ws.Cells.Item(103, 10).Value = 1749.5  # J103: 2499.6667 -> 1749.5
$ws.Cells.Item(103, 11).Value = 3000  # K103: 11250 -> 3000
$ws.Cells.Item(103, 12).Value = 5248.5  # L103: 7499.000100000001 -> 5248.5
$ws.Cells.Item(103, 13).Value = -2414  # M103: -10664 -> -2414
$ws.Cells.Item(103, 14).Value = -6420.5  # N103: -8671.000100000001 -> -6420.5

# Row 113 (ALC)
$ws.Cells.Item(113, 8).Value = 2605  # H113: 6854.2 -> 2605
$ws.Cells.Item(113, 9).Value = 2605  # I113: 8201.666999999999 -> 2605
$ws.Cells.Item(113, 10).Value = 0  # J113: 4833 -> 0
$ws.Cells.Item(113, 11).Value = 2605  # K113: 8201.666999999999 -> 2605
$ws.Cells.Item(113, 12).Value = 0  # L113: 4833 -> 0
$ws.Cells.Item(113, 13).Value = 649  # M113: -4947.666999999999 -> 649
$ws.Cells.Item(113, 14).ClearContents()  # N113: -11341 -> (removed)

# Row 132 (ALC)
$ws.Cells.Item(132, 8).Value = 5535.9165  # H132: 6470.1 -> 5535.9165
$ws.Cells.Item(132, 9).Value = 1643.1  # I132: 1837.625 -> 1643.1
$ws.Cells.Item(132, 11).Value = 4929.299999999999  # K132: 5512.875 -> 4929.299999999999
$ws.Cells.Item(132, 13).Value = -2399.299999999999  # M132: -2982.875 -> -2399.299999999999

# Row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 3077.2354  # H137: 2920.8 -> 3077.2354
$ws.Cells.Item(137, 9).Value = 1718.5  # I137: 1840.2 -> 1718.5
$ws.Cells.Item(137, 10).Value = 3818.3635  # J137: 3281 -> 3818.3635
$ws.Cells.Item(137, 11).Value = 5155.5  # K137: 5520.6 -> 5155.5
$ws.Cells.Item(137, 12).Value = 11455.0905  # L137: 9843 -> 11455.0905
$ws.Cells.Item(137, 13).Value = -2605.5  # M137: -2970.6 -> -2605.5
$ws.Cells.Item(137, 14).Value = -16555.0905  # N137: -14943 -> -16555.0905

$ws = $wb.Worksheets.Item("ARM")
# Row 16 (ARM)
$ws.Cells.Item(16, 8).Value = 285.75  # H16: 406 -> 285.75
$ws.Cells.Item(16, 9).Value = 219  # I16: 406 -> 219
$ws.Cells.Item(16, 10).Value = 352.5  # J16: 0 -> 352.5
$ws.Cells.Item(16, 11).Value = 219  # K16: 406 -> 219
$ws.Cells.Item(16, 12).Value = 352.5  # L16: 0 -> 352.5
$ws.Cells.Item(16, 13).Value = 68  # M16: -119 -> 68
$ws.Cells.Item(16, 14).Value = -926.5  # N16: None -> -926.5

# Row 135 (ARM)
$ws.Cells.Item(135, 8).Value = 100429  # H135: 60000 -> 100429
$ws.Cells.Item(135, 10).Value = 100429  # J135: 60000 -> 100429
$ws.Cells.Item(135, 12).Value = 100429  # L135: 60000 -> 100429
$ws.Cells.Item(135, 14).Value = -110569  # N135: -70140 -> -110569

$ws = $wb.Worksheets.Item("CRP")
# Row 5 (CRP)
$ws.Cells.Item(5, 8).Value = 1065.8  # H5: 851 -> 1065.8
$ws.Cells.Item(5, 10).Value = 2999.3333  # J5: 2999.5 -> 2999.3333
$ws.Cells.Item(5, 12).Value = 2999.3333  # L5: 2999.5 -> 2999.3333
$ws.Cells.Item(5, 14).Value = -3223.3333  # N5: -3223.5 -> -3223.3333

# Row 16 (CRP)
$ws.Cells.Item(16, 8).Value = 1272.9412  # H16: 1425.5385 -> 1272.9412
$ws.Cells.Item(16, 9).Value = 1185.8182  # I16: 1303.6666 -> 1185.8182
$ws.Cells.Item(16, 10).Value = 1432.6666  # J16: 1699.75 -> 1432.6666
$ws.Cells.Item(16, 11).Value = 1185.8182  # K16: 1303.6666 -> 1185.8182
$ws.Cells.Item(16, 12).Value = 1432.6666  # L16: 1699.75 -> 1432.6666
$ws.Cells.Item(16, 13).Value = -898.8181999999999  # M16: -1016.6666 -> -898.8181999999999
$ws.Cells.Item(16, 14).Value = -2006.6666  # N16: -2273.75 -> -2006.6666

# Row 31 (CRP)
$ws.Cells.Item(31, 8).Value = 8254.280000000001  # H31: 7902.1113 -> 8254.280000000001
$ws.Cells.Item(31, 9).Value = 5068  # I31: 4440.8 -> 5068
$ws.Cells.Item(31, 11).Value = 5068  # K31: 4440.8 -> 5068
$ws.Cells.Item(31, 13).Value = -4773  # M31: -4145.8 -> -4773

# Row 34 (CRP)
$ws.Cells.Item(34, 8).Value = 8254.280000000001  # H34: 7902.1113 -> 8254.280000000001
$ws.Cells.Item(34, 9).Value = 5068  # I34: 4440.8 -> 5068
$ws.Cells.Item(34, 11).Value = 5068  # K34: 4440.8 -> 5068
$ws.Cells.Item(34, 13).Value = -4866  # M34: -4238.8 -> -4866

# Row 100 (CRP)
$ws.Cells.Item(100, 8).Value = 79999.5  # H100: 80000 -> 79999.5
$ws.Cells.Item(100, 10).Value = 79999.5  # J100: 80000 -> 79999.5
$ws.Cells.Item(100, 12).Value = 79999.5  # L100: 80000 -> 79999.5
$ws.Cells.Item(100, 14).Value = -82163.5  # N100: -82164 -> -82163.5

# Row 108 (CRP)
$ws.Cells.Item(108, 8).Value = 80210.39999999999  # H108: 67541.664 -> 80210.39999999999
$ws.Cells.Item(108, 9).Value = 0  # I108: 22625 -> 0
$ws.Cells.Item(108, 10).Value = 80210.39999999999  # J108: 90000 -> 80210.39999999999
$ws.Cells.Item(108, 11).Value = 0  # K108: 22625 -> 0
$ws.Cells.Item(108, 12).Value = 80210.39999999999  # L108: 90000 -> 80210.39999999999
$ws.Cells.Item(108, 13).ClearContents()  # M108: -18785 -> (removed)
$ws.Cells.Item(108, 14).Value = -87890.39999999999  # N108: -97680 -> -87890.39999999999

# Row 113 (CRP)
$ws.Cells.Item(113, 8).Value = 1272.9412  # H113: 1425.5385 -> 1272.9412
$ws.Cells.Item(113, 9).Value = 1185.8182  # I113: 1303.6666 -> 1185.8182
$ws.Cells.Item(113, 10).Value = 1432.6666  # J113: 1699.75 -> 1432.6666
$ws.Cells.Item(113, 11).Value = 1185.8182  # K113: 1303.6666 -> 1185.8182
$ws.Cells.Item(113, 12).Value = 1432.6666  # L113: 1699.75 -> 1432.6666
$ws.Cells.Item(113, 13).Value = 984.1818000000001  # M113: 866.3334 -> 984.1818000000001
$ws.Cells.Item(113, 14).Value = -5772.6666  # N113: -6039.75 -> -5772.6666

# Row 132 (CRP)
$ws.Cells.Item(132, 8).Value = 4688.5  # H132: 4749.154 -> 4688.5
$ws.Cells.Item(132, 9).Value = 4089.5  # I132: 4110.5557 -> 4089.5
$ws.Cells.Item(132, 11).Value = 12268.5  # K132: 12331.6671 -> 12268.5
$ws.Cells.Item(132, 13).Value = -9738.5  # M132: -9801.667099999999 -> -9738.5

# Row 134 (CRP)
$ws.Cells.Item(134, 8).Value = 1666.3334  # H134: 1861.125 -> 1666.3334
$ws.Cells.Item(134, 9).Value = 1666.3334  # I134: 1861.125 -> 1666.3334
$ws.Cells.Item(134, 11).Value = 4999.0002  # K134: 5583.375 -> 4999.0002
$ws.Cells.Item(134, 13).Value = -2464.0002  # M134: -3048.375 -> -2464.0002

$ws = $wb.Worksheets.Item("CUL")
# Row 40 (CUL)
$ws.Cells.Item(40, 8).Value = 113.4  # H40: 92.5 -> 113.4
$ws.Cells.Item(40, 10).Value = 197.5  # J40: 198 -> 197.5
$ws.Cells.Item(40, 12).Value = 790  # L40: 792 -> 790
$ws.Cells.Item(40, 14).Value = -928  # N40: -930 -> -928

# Row 60 (CUL)
$ws.Cells.Item(60, 8).Value = 1135.3334  # H60: 1231.25 -> 1135.3334
$ws.Cells.Item(60, 10).Value = 2766  # J60: 2750 -> 2766
$ws.Cells.Item(60, 12).Value = 8298  # L60: 8250 -> 8298
$ws.Cells.Item(60, 14).Value = -8800  # N60: -8752 -> -8800

# Row 80 (CUL)
$ws.Cells.Item(80, 8).Value = 4947.913  # H80: 4984 -> 4947.913
$ws.Cells.Item(80, 10).Value = 5125.75  # J80: 5333.25 -> 5125.75
$ws.Cells.Item(80, 12).Value = 15377.25  # L80: 15999.75 -> 15377.25
$ws.Cells.Item(80, 14).Value = -17249.25  # N80: -17871.75 -> -17249.25

# Row 83 (CUL)
$ws.Cells.Item(83, 8).Value = 4947.913  # H83: 4984 -> 4947.913
$ws.Cells.Item(83, 10).Value = 5125.75  # J83: 5333.25 -> 5125.75
$ws.Cells.Item(83, 12).Value = 46131.75  # L83: 47999.25 -> 46131.75
$ws.Cells.Item(83, 14).Value = -55491.75  # N83: -57359.25 -> -55491.75

# Row 128 (CUL)
$ws.Cells.Item(128, 8).Value = 533328.2  # H128: 521418 -> 533328.2
$ws.Cells.Item(128, 9).Value = 533328.2  # I128: 521418 -> 533328.2
$ws.Cells.Item(128, 11).Value = 1599984.6  # K128: 1564254 -> 1599984.6
$ws.Cells.Item(128, 13).Value = -1595004.6  # M128: -1559274 -> -1595004.6

$ws = $wb.Worksheets.Item("GSM")
# Row 6 (GSM)
$ws.Cells.Item(6, 8).Value = 182.5  # H6: 154.5 -> 182.5
$ws.Cells.Item(6, 10).Value = 182.5  # J6: 154.5 -> 182.5
$ws.Cells.Item(6, 12).Value = 182.5  # L6: 154.5 -> 182.5
$ws.Cells.Item(6, 14).Value = -408.5  # N6: -380.5 -> -408.5

# Row 9 (GSM)
$ws.Cells.Item(9, 8).Value = 221.11111  # H9: 230.5 -> 221.11111
$ws.Cells.Item(9, 9).Value = 230  # I9: 230.5 -> 230
$ws.Cells.Item(9, 10).Value = 150  # J9: 0 -> 150
$ws.Cells.Item(9, 11).Value = 230  # K9: 230.5 -> 230
$ws.Cells.Item(9, 12).Value = 150  # L9: 0 -> 150
$ws.Cells.Item(9, 13).Value = -60  # M9: -60.5 -> -60
$ws.Cells.Item(9, 14).Value = -490  # N9: None -> -490

# Row 16 (GSM)
$ws.Cells.Item(16, 8).Value = 182.5  # H16: 154.5 -> 182.5
$ws.Cells.Item(16, 10).Value = 182.5  # J16: 154.5 -> 182.5
$ws.Cells.Item(16, 12).Value = 182.5  # L16: 154.5 -> 182.5
$ws.Cells.Item(16, 14).Value = -682.5  # N16: -654.5 -> -682.5

# Row 23 (GSM)
$ws.Cells.Item(23, 8).Value = 822.5  # H23: 951 -> 822.5
$ws.Cells.Item(23, 10).Value = 822.5  # J23: 951 -> 822.5
$ws.Cells.Item(23, 12).Value = 822.5  # L23: 951 -> 822.5
$ws.Cells.Item(23, 14).Value = -1268.5  # N23: -1397 -> -1268.5

# Row 80 (GSM)
$ws.Cells.Item(80, 8).Value = 999.5  # H80: 0 -> 999.5
$ws.Cells.Item(80, 9).Value = 999  # I80: 0 -> 999
$ws.Cells.Item(80, 10).Value = 1000  # J80: 0 -> 1000
$ws.Cells.Item(80, 11).Value = 999  # K80: 0 -> 999
$ws.Cells.Item(80, 12).Value = 1000  # L80: 0 -> 1000
$ws.Cells.Item(80, 13).Value = -1  # M80: None -> -1
$ws.Cells.Item(80, 14).Value = -2996  # N80: None -> -2996

# Row 83 (GSM)
$ws.Cells.Item(83, 8).Value = 999.5  # H83: 0 -> 999.5
$ws.Cells.Item(83, 9).Value = 999  # I83: 0 -> 999
$ws.Cells.Item(83, 10).Value = 1000  # J83: 0 -> 1000
$ws.Cells.Item(83, 11).Value = 4995  # K83: 0 -> 4995
$ws.Cells.Item(83, 12).Value = 5000  # L83: 0 -> 5000
$ws.Cells.Item(83, 13).Value = -3  # M83: None -> -3
$ws.Cells.Item(83, 14).Value = -14984  # N83: None -> -14984

# Row 122 (GSM)
$ws.Cells.Item(122, 8).Value = 2628.2222  # H122: 2678.3333 -> 2628.2222
$ws.Cells.Item(122, 9).Value = 665.25  # I122: 778 -> 665.25
$ws.Cells.Item(122, 11).Value = 1995.75  # K122: 2334 -> 1995.75
$ws.Cells.Item(122, 13).Value = 454.25  # M122: 116 -> 454.25

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Cells.Item(7, 8).Value = 4764.1333  # H7: 4764.8 -> 4764.1333
$ws.Cells.Item(7, 9).Value = 3133.875  # I7: 3135.125 -> 3133.875
$ws.Cells.Item(7, 11).Value = 3133.875  # K7: 3135.125 -> 3133.875
$ws.Cells.Item(7, 13).Value = -3021.875  # M7: -3023.125 -> -3021.875

# Row 93 (LTW)
$ws.Cells.Item(93, 8).Value = 1976  # H93: 2301.5 -> 1976
$ws.Cells.Item(93, 9).Value = 1953.7142  # I93: 2499 -> 1953.7142
$ws.Cells.Item(93, 10).Value = 2054  # J93: 2104 -> 2054
$ws.Cells.Item(93, 11).Value = 1953.7142  # K93: 2499 -> 1953.7142
$ws.Cells.Item(93, 12).Value = 2054  # L93: 2104 -> 2054
$ws.Cells.Item(93, 13).Value = -705.7141999999999  # M93: -1251 -> -705.7141999999999
$ws.Cells.Item(93, 14).Value = -4550  # N93: -4600 -> -4550

# Row 126 (LTW)
$ws.Cells.Item(126, 8).Value = 4764.1333  # H126: 4764.8 -> 4764.1333
$ws.Cells.Item(126, 9).Value = 3133.875  # I126: 3135.125 -> 3133.875
$ws.Cells.Item(126, 11).Value = 9401.625  # K126: 9405.375 -> 9401.625
$ws.Cells.Item(126, 13).Value = -6931.625  # M126: -6935.375 -> -6931.625

$ws = $wb.Worksheets.Item("WVR")
# Row 42 (WVR)
$ws.Cells.Item(42, 8).Value = 32500  # H42: 0 -> 32500
$ws.Cells.Item(42, 9).Value = 32500  # I42: 0 -> 32500
$ws.Cells.Item(42, 10).Value = 32500  # J42: 0 -> 32500
$ws.Cells.Item(42, 11).Value = 32500  # K42: 0 -> 32500
$ws.Cells.Item(42, 12).Value = 32500  # L42: 0 -> 32500
$ws.Cells.Item(42, 13).Value = -32122  # M42: None -> -32122
$ws.Cells.Item(42, 14).Value = -33256  # N42: None -> -33256

# Row 81 (WVR)
$ws.Cells.Item(81, 8).Value = 0  # H81: 3499 -> 0
$ws.Cells.Item(81, 10).Value = 0  # J81: 3499 -> 0
$ws.Cells.Item(81, 12).Value = 0  # L81: 6998 -> 0
$ws.Cells.Item(81, 14).ClearContents()  # N81: -9120 -> (removed)

# Row 84 (WVR)
$ws.Cells.Item(84, 8).Value = 0  # H84: 3499 -> 0
$ws.Cells.Item(84, 10).Value = 0  # J84: 3499 -> 0
$ws.Cells.Item(84, 12).Value = 0  # L84: 34990 -> 0
$ws.Cells.Item(84, 14).ClearContents()  # N84: -45598 -> (removed)

# Row 107 (WVR)
$ws.Cells.Item(107, 8).Value = 691.1667  # H107: 849.5 -> 691.1667
$ws.Cells.Item(107, 9).Value = 788.2222  # I107: 849.5 -> 788.2222
$ws.Cells.Item(107, 10).Value = 400  # J107: 0 -> 400
$ws.Cells.Item(107, 11).Value = 2364.6666  # K107: 2548.5 -> 2364.6666
$ws.Cells.Item(107, 12).Value = 1200  # L107: 0 -> 1200
$ws.Cells.Item(107, 13).Value = -444.6666  # M107: -628.5 -> -444.6666
$ws.Cells.Item(107, 14).Value = -5040  # N107: None -> -5040
